# Weekly update: insert a new price record for "Albahaca" (Terminal La Palmera
# de La Serena) as row 50, shifting the existing historical rows (50-97) down
# by one (to 51-98). The new row reuses the same market / product attributes
# as the previous entry at that position, with refreshed date and prices.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 50; everything below (50:97) moves to (51:98)
$ws.Rows("50:50").Insert()

# Fill in the new weekly record in row 50
$ws.Cells.Item(50, 1).Value2  = 8
$ws.Cells.Item(50, 2).Value   = "Terminal La Palmera de La Serena"
$ws.Cells.Item(50, 3).Value   = "Coquimbo"
$ws.Cells.Item(50, 4).Value2  = 44741
$ws.Cells.Item(50, 5).Value2  = 4
$ws.Cells.Item(50, 6).Value2  = 100112052
$ws.Cells.Item(50, 7).Value   = "Albahaca"
$ws.Cells.Item(50, 8).Value   = "Sin especificar"
$ws.Cells.Item(50, 9).Value   = "Primera"
$ws.Cells.Item(50, 10).Value2 = 800
$ws.Cells.Item(50, 11).Value2 = 3500
$ws.Cells.Item(50, 12).Value2 = 4000
$ws.Cells.Item(50, 13).Value2 = 3750
$ws.Cells.Item(50, 14).Value  = "`$/paquete"
$ws.Cells.Item(50, 15).Value  = "Región de Arica y Parinacota"
$ws.Cells.Item(50, 16).Value2 = 3750
$ws.Cells.Item(50, 17).Value2 = 1
$ws.Cells.Item(50, 18).Value  = "Hortaliza"
